$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update aggregate metrics to reflect the new closed trade.
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.98   # Current Capital
$summary.Range("B4").Value = -0.02     # Total P&L $
$summary.Range("B5").Value = -0.03     # Total P&L %
$summary.Range("B6").Value = 12        # Total Trades
$summary.Range("B7").Value = 5         # Winning Trades
$summary.Range("B9").Value = 41.67     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4).
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.98      # Capital
$status.Range("D4").Value = 12         # Trades
$status.Range("E4").Value = -0.02      # P&L $
$status.Range("F4").Value = -0.02      # P&L %
$status.Range("G4").Value = 41.67      # Win Rate %

# ---------------------------------------------------------------------------
# Sheets "All Trades" and "MarketMaking": append the newly closed trade #12
# as row 13 on both sheets (they mirror the same trade log).
# ---------------------------------------------------------------------------
$sheetNames = @("All Trades", "MarketMaking")
foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    $ws.Range("A13").Value = 12

    # Force the date/time-looking text to stay as literal text instead of
    # being auto-converted to a date serial number.
    $ws.Range("B13").NumberFormat = "@"
    $ws.Range("B13").Value = "2026-02-17"

    $ws.Range("C13").Value = "07:58:35"
    $ws.Range("D13").Value = "MarketMaking"
    $ws.Range("E13").Value = "DOWN"
    $ws.Range("F13").Value = 0.93
    $ws.Range("G13").Value = 0.95
    $ws.Range("H13").Value = "CLOSED"
    $ws.Range("I13").Value = 2.1505
    $ws.Range("J13").Value = 0.02
    $ws.Range("K13").Value = 99.98
    $ws.Range("L13").Value = 0
    $ws.Range("M13").Value = 0
    $ws.Range("N13").Value = 0.6
    $ws.Range("O13").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P13").Value = "early_exit"
    $ws.Range("Q13").Value = 0.13
}
